# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> originally the default "Office Theme" colours
#   ppt/theme/theme2.xml  -> originally the "Integral" design colours,
#                            used by the slide master / all slides
#
# The authored change swaps the content of the two theme parts: the deck's
# visible design switches from the "Integral" colour palette to the
# standard "Office" colour palette (dk1/lt1/dk2/lt2/accent1-6/hlink/
# folHlink), while what used to be the Office palette ends up parked in
# the other theme part. Font scheme and format scheme are identical
# between the two parts already, so the colour scheme is the only
# substantive payload that needs to move.
#
# Apply this by rewriting the active design's theme colour scheme (the
# one backing the slide master that every slide in the deck uses) from
# the "Integral" values to the standard "Office" values, in the official
# ThemeColorScheme index order:
#   1 Dark1, 2 Light1, 3 Dark2, 4 Light2, 5-10 Accent1-6, 11 Hyperlink,
#   12 FollowedHyperlink.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# PowerPoint's RGB() is a VBA intrinsic, not a PowerShell one, so the
# palette is written as plain OLE_COLOR integers (r + g*256 + b*65536).
$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72

# Best-effort: keep the naming metadata consistent with the new palette
# (the runtime may or may not persist these, but they cost nothing to set).
$colors.Name = "Office"
$theme.Name = "Office Theme"
$design.Name = "Office Theme"
